$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 217, pushing the existing rows 217-231 down to 218-232.
$ws.Rows.Item(217).Insert()

# Fill in the new row 217 with a new weekly price record (same market/category
# metadata as its neighbours, new date + volume + price figures).
$ws.Range("A217").Value = 4
$ws.Range("B217").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C217").Value = "Los Lagos"
$ws.Range("D217").Value = 44746
$ws.Range("E217").Value = 10
$ws.Range("F217").Value = 100112039
$ws.Range("G217").Value = "Ciboulette"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 40
$ws.Range("K217").Value = 2500
$ws.Range("L217").Value = 2500
$ws.Range("M217").Value = 2500
$ws.Range("N217").Value = "`$/docena de atados"
$ws.Range("O217").Value = "Región Metropolitana"
$ws.Range("P217").Value = 833
$ws.Range("Q217").Value = 3
$ws.Range("R217").Value = "Hortaliza"
